# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Sheets workbook per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 138.53847
$ws.Range("I9").Value = 149.18182
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 149.18182
$ws.Range("L9").Value = 80
$ws.Range("M9").Value = 19.81818000000001
$ws.Range("N9").Value = -418
$ws.Range("H17").Value = 403438.62
$ws.Range("J17").Value = 403438.62
$ws.Range("L17").Value = 1210315.86
$ws.Range("N17").Value = -1210651.86
$ws.Range("H33").Value = 170.25
$ws.Range("I33").Value = 175.84616
$ws.Range("J33").Value = 146
$ws.Range("K33").Value = 175.84616
$ws.Range("L33").Value = 146
$ws.Range("M33").Value = 53.15384
$ws.Range("N33").Value = -604
$ws.Range("H100").Value = 6149.7666
$ws.Range("I100").Value = 832.6111
$ws.Range("J100").Value = 14125.5
$ws.Range("K100").Value = 832.6111
$ws.Range("L100").Value = 14125.5
$ws.Range("M100").Value = -291.6111
$ws.Range("N100").Value = -15207.5
$ws.Range("H112").Value = 1003111.6
$ws.Range("J112").Value = 1003111.6
$ws.Range("L112").Value = 3009334.8
$ws.Range("N112").Value = -3011550.8
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120
$ws.Range("H135").Value = 10185.357
$ws.Range("I135").Value = 9219.666999999999
$ws.Range("K135").Value = 82977.003
$ws.Range("M135").Value = -80442.003
$ws.Range("H137").Value = 2609.238
$ws.Range("I137").Value = 2544.111
$ws.Range("K137").Value = 7632.333
$ws.Range("M137").Value = -5082.333
$ws.Range("H138").Value = 345639
$ws.Range("I138").Value = 43083.293
$ws.Range("J138").Value = 1253306.1
$ws.Range("K138").Value = 129249.879
$ws.Range("L138").Value = 3759918.3
$ws.Range("M138").Value = -124109.879
$ws.Range("N138").Value = -3770198.3
$ws.Range("H141").Value = 2844.32
$ws.Range("J141").Value = 5935.625
$ws.Range("L141").Value = 17806.875
$ws.Range("N141").Value = -28166.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = -4626
$ws.Range("H32").Value = 7280.15
$ws.Range("I32").Value = 6703.4062
$ws.Range("J32").Value = 9587.125
$ws.Range("K32").Value = 6703.4062
$ws.Range("L32").Value = 9587.125
$ws.Range("M32").Value = -6416.4062
$ws.Range("N32").Value = -10161.125
$ws.Range("H38").Value = 9666.666999999999
$ws.Range("I38").Value = 9666.666999999999
$ws.Range("K38").Value = 9666.666999999999
$ws.Range("M38").Value = -9199.666999999999
$ws.Range("H131").Value = 75446.82000000001
$ws.Range("J131").Value = 75446.82000000001
$ws.Range("L131").Value = 75446.82000000001
$ws.Range("N131").Value = -85526.82000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4158.579
$ws.Range("I86").Value = 5042.2856
$ws.Range("J86").Value = 1684.2
$ws.Range("K86").Value = 5042.2856
$ws.Range("L86").Value = 1684.2
$ws.Range("M86").Value = -3919.2856
$ws.Range("N86").Value = -3930.2
$ws.Range("H89").Value = 4158.579
$ws.Range("I89").Value = 5042.2856
$ws.Range("J89").Value = 1684.2
$ws.Range("K89").Value = 25211.428
$ws.Range("L89").Value = 8421
$ws.Range("M89").Value = -19595.428
$ws.Range("N89").Value = -19653
$ws.Range("H94").Value = 3242.6428
$ws.Range("J94").Value = 4641
$ws.Range("L94").Value = 4641
$ws.Range("N94").Value = -5543
$ws.Range("H117").Value = 143870.5
$ws.Range("J117").Value = 143870.5
$ws.Range("L117").Value = 143870.5
$ws.Range("N117").Value = -153048.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 8305.6
$ws.Range("I2").Value = 8305.6
$ws.Range("K2").Value = 8305.6
$ws.Range("M2").Value = -8192.6
$ws.Range("H36").Value = 28308.25
$ws.Range("I36").Value = 13000
$ws.Range("J36").Value = 74233
$ws.Range("K36").Value = 13000
$ws.Range("L36").Value = 74233
$ws.Range("M36").Value = -12612
$ws.Range("N36").Value = -75009
$ws.Range("H40").Value = 28308.25
$ws.Range("I40").Value = 13000
$ws.Range("J40").Value = 74233
$ws.Range("K40").Value = 13000
$ws.Range("L40").Value = 74233
$ws.Range("M40").Value = -12840
$ws.Range("N40").Value = -74553
$ws.Range("H134").Value = 2192.2363
$ws.Range("I134").Value = 1217.881
$ws.Range("K134").Value = 3653.643
$ws.Range("M134").Value = -1118.643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11100824
$ws.Range("I4").Value = 11782973
$ws.Range("K4").Value = 35348919
$ws.Range("M4").Value = -35348807
$ws.Range("H5").Value = 4489.241
$ws.Range("I5").Value = 378
$ws.Range("J5").Value = 9549.23
$ws.Range("K5").Value = 1134
$ws.Range("L5").Value = 28647.69
$ws.Range("M5").Value = -1022
$ws.Range("N5").Value = -28871.69
$ws.Range("H9").Value = 533731.6
$ws.Range("I9").Value = 1000000
$ws.Range("J9").Value = 222886
$ws.Range("K9").Value = 3000000
$ws.Range("L9").Value = 668658
$ws.Range("M9").Value = -2999776
$ws.Range("N9").Value = -669106
$ws.Range("H40").Value = 68.55556
$ws.Range("I40").Value = 68.55556
$ws.Range("K40").Value = 274.22224
$ws.Range("M40").Value = -205.22224
$ws.Range("H69").Value = 4011.5
$ws.Range("I69").Value = 4011.5
$ws.Range("K69").Value = 12034.5
$ws.Range("M69").Value = -11223.5
$ws.Range("H72").Value = 4011.5
$ws.Range("I72").Value = 4011.5
$ws.Range("K72").Value = 36103.5
$ws.Range("M72").Value = -32047.5
$ws.Range("H121").Value = 4538.7
$ws.Range("J121").Value = 4959.6665
$ws.Range("L121").Value = 14878.9995
$ws.Range("N121").Value = -17498.9995
$ws.Range("H134").Value = 2651.4
$ws.Range("I134").Value = 2168.2222
$ws.Range("K134").Value = 6504.6666
$ws.Range("M134").Value = -1434.6666
$ws.Range("H135").Value = 4489.241
$ws.Range("I135").Value = 378
$ws.Range("J135").Value = 9549.23
$ws.Range("K135").Value = 3402
$ws.Range("L135").Value = 85943.06999999999
$ws.Range("M135").Value = -867
$ws.Range("N135").Value = -91013.06999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 6670368
$ws.Range("I14").Value = 10005000
$ws.Range("J14").Value = 1105
$ws.Range("K14").Value = 10005000
$ws.Range("L14").Value = 1105
$ws.Range("M14").Value = -10004832
$ws.Range("N14").Value = -1441
$ws.Range("H46").Value = 48849
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H57").Value = 68648.94
$ws.Range("J57").Value = 92125.73
$ws.Range("L57").Value = 92125.73
$ws.Range("N57").Value = -93765.73
$ws.Range("H117").Value = 59247.75
$ws.Range("J117").Value = 59247.75
$ws.Range("L117").Value = 59247.75
$ws.Range("N117").Value = -66131.75
$ws.Range("H132").Value = 7755042
$ws.Range("I132").Value = 9012119
$ws.Range("K132").Value = 27036357
$ws.Range("M132").Value = -27033827
$ws.Range("H134").Value = 88500
$ws.Range("J134").Value = 88500
$ws.Range("L134").Value = 265500
$ws.Range("N134").Value = -270570
$ws.Range("H141").Value = 135214
$ws.Range("J141").Value = 135214
$ws.Range("L141").Value = 135214
$ws.Range("N141").Value = -145574

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 232.5
$ws.Range("I55").Value = 132.28572
$ws.Range("J55").Value = 332.7143
$ws.Range("K55").Value = 132.28572
$ws.Range("L55").Value = 332.7143
$ws.Range("M55").Value = 40.71428
$ws.Range("N55").Value = -678.7143
$ws.Range("H114").Value = 89627.336
$ws.Range("J114").Value = 89627.336
$ws.Range("L114").Value = 89627.336
$ws.Range("N114").Value = -98305.336
$ws.Range("H122").Value = 2867.9697
$ws.Range("I122").Value = 2298.7407
$ws.Range("K122").Value = 6896.222099999999
$ws.Range("M122").Value = -4446.222099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 5333
$ws.Range("I6").Value = 499.5
$ws.Range("J6").Value = 15000
$ws.Range("K6").Value = 499.5
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = -384.5
$ws.Range("N6").Value = -15230
$ws.Range("H62").Value = 205444.25
$ws.Range("I62").Value = 205444.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 205444.25
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -204820.25
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 205444.25
$ws.Range("I65").Value = 205444.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1027221.25
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -1024101.25
$ws.Range("N65").ClearContents()
$ws.Range("H101").Value = 63813.5
$ws.Range("J101").Value = 63813.5
$ws.Range("L101").Value = 63813.5
$ws.Range("N101").Value = -70303.5
$ws.Range("H113").Value = 5953315
$ws.Range("J113").Value = 1198.125
$ws.Range("L113").Value = 3594.375
$ws.Range("N113").Value = -7934.375
